$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.523.74"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.696.53"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'316.74"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.3907"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").Value = "'0.4080"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").Value = "'1.492"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").Value = "'1.002"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").Value = "'53.20"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'0.08837"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'26.56"
$ws.Range("E13").Value = "  +12.44%  "
$ws.Range("D14").Value = "'7.488"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "'8.321"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "'0.00001370"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("D17").Value = "1.695.13"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "'98.25"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").Value = "'0.07214"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").Value = "'20.57"
$ws.Range("E20").Value = "  +4.60%  "
$ws.Range("D21").Value = "'7.323"
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "'14.39"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").Value = "24.522.58"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "'3.024"
$ws.Range("E25").Value = "  -3.55%  "
$ws.Range("D26").Value = "'2.336"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").Value = "'23.05"
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("D28").Value = "'168.49"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("D29").Value = "'147.34"
$ws.Range("E29").Value = "  +8.66%  "
$ws.Range("D30").Value = "'8.497"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("D31").Value = "'5.376"
$ws.Range("E31").Value = "  +4.38%  "
$ws.Range("D32").Value = "1.880.60"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.08821"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.179"
$ws.Range("E34").Value = "  +10.48%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'7.265"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.052"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").Value = "'0.03078"
$ws.Range("E37").Value = "  +10.61%  "
$ws.Range("D38").Value = "'0.2810"
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("D39").Value = "'10.96"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").Value = "'0.09192"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").Value = "'14.28"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "'0.8007"
$ws.Range("E42").Value = "  +4.39%  "
$ws.Range("D43").Value = "'1.483"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "'17.63"
$ws.Range("E44").Value = "  +11.30%  "
$ws.Range("D45").Value = "'2.686"
$ws.Range("E45").Value = "  +5.35%  "
$ws.Range("D46").Value = "'0.7269"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").Value = "'4.274"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("D48").Value = "'1.406"
$ws.Range("E48").Value = "  +4.50%  "
$ws.Range("D49").Value = "'1.001"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "'140.66"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").Value = "'0.08207"
$ws.Range("E51").Value = "  +2.92%  "
